$wb = $excel.ActiveWorkbook

$wsNotes = $wb.Worksheets.Item("Data Deliverable Notes")

# Rows whose "Completed?" (column C) flag changes from "x" to "AM" (initials)
$rows = @(17, 18, 22, 23, 29, 34, 35, 36, 37, 38, 39, 40, 41)
foreach ($r in $rows) {
    $wsNotes.Range("C$r").Value = "AM"
}

# Make "Data Deliverable Notes" the active (selected) sheet in the workbook,
# with C41 as the last-selected cell (tracks where the final edit was made)
$wsNotes.Activate()
$wsNotes.Range("C41").Select()
